# Update the "取得日時" (acquisition timestamp) column on the first sheet
# ("ランサーズ") so that rows 2-8 reflect the new scrape timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-08 01:22:52"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
